# Auto-generated edit script applying the cryptos.xlsx diff
# (price refresh + two coin-row reorderings) via Excel COM interop.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as literal text in the source file (e.g. "0.110",
# "76.468.38" using "." as a thousands separator). Excel auto-converts any
# plain numeric-looking string typed into a General cell to a Number, which
# would silently drop things like trailing zeros ("0.110" -> 0.11). To keep
# those specific cells as real text (matching the workbook author's data),
# we flip them to the Text number format before writing the new value.

# Row 2
$ws.Range("D2").Value = "76.526.87"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "2.929.51"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.68"
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "594.45"
$ws.Range("E6").Value = "  -1.19%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("E8").Value = "  -1.60%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.198"
$ws.Range("E9").Value = "  +2.43%  "

# Row 10
$ws.Range("D10").Value = "2.925.76"
$ws.Range("E10").Value = "  +0.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.440"
$ws.Range("E11").Value = "  +10.85%  "

# Row 12
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("E13").Value = "  -1.43%  "

# Row 14
$ws.Range("D14").Value = "3.467.31"
$ws.Range("E14").Value = "  +0.94%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.35"
$ws.Range("E15").Value = "  +2.53%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "76.423.51"
$ws.Range("E16").Value = "  +0.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000189"
$ws.Range("E17").Value = "  -1.52%  "

# Row 18
$ws.Range("D18").Value = "2.912.52"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.46"
$ws.Range("E19").Value = "  +6.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.73"
$ws.Range("E20").Value = "  -3.05%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.57"
$ws.Range("E21").Value = "  -2.85%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  +3.49%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.25"
$ws.Range("E23").Value = "  -2.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.82"
$ws.Range("E24").Value = "  -0.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").Value = "3.077.44"
$ws.Range("E26").Value = "  +1.48%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.26"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.62"
$ws.Range("E28").Value = "  -2.38%  "

# Row 29
$ws.Range("E29").Value = "  -2.14%  "

# Row 30
$ws.Range("E30").Value = "  +0.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +6.08%  "

# Row 32
$ws.Range("E32").Value = "  -3.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "499.73"
$ws.Range("E33").Value = "  -3.03%  "

# Row 34
$ws.Range("E34").Value = "  -0.25%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.16%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.31"
$ws.Range("E36").Value = "  -0.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.10"
$ws.Range("E37").Value = "  -0.78%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.391"
$ws.Range("E38").Value = "  +12.03%  "

# Row 39
$ws.Range("B39").Value = "Cronos"
$ws.Range("C39").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.110"
$ws.Range("E39").Value = "  +18.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.95"
$ws.Range("E40").Value = "  +1.36%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  -4.91%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "178.57"
$ws.Range("E43").Value = "  -3.02%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.91"
$ws.Range("E44").Value = "  -3.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.65"
$ws.Range("E45").Value = "  -2.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.94"
$ws.Range("E46").Value = "  -0.78%  "

# Row 47
$ws.Range("E47").Value = "  -4.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.591"
$ws.Range("E48").Value = "  +0.43%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -2.89%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.85"
$ws.Range("E50").Value = "  +1.83%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.658"
$ws.Range("E51").Value = "  -3.94%  "
